$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet had three accidental duplicate track rows; remove them.
# Deleting from the bottom up keeps the remaining row numbers stable
# while we work.
$ws.Rows(111).Delete()
$ws.Rows(107).Delete()
$ws.Rows(80).Delete()

# Populate the "Year out" values (column E) for every data row that was
# still missing one, using the final (post-delete) row numbers.
$yearData = @{
    56 = 2018
    57 = 2015
    58 = 2017
    59 = 2015
    60 = 2018
    61 = 2019
    62 = 2020
    63 = 2018
    64 = 2017
    65 = 2019
    66 = 2014
    67 = 2013
    68 = 2011
    69 = 2020
    70 = 2017
    71 = 2012
    72 = 2015
    73 = 2018
    74 = 2014
    75 = 2014
    76 = 2019
    77 = 2020
    78 = 2016
    79 = 2019
    80 = 2015
    81 = 2020
    82 = 2014
    83 = 2020
    84 = 2015
    85 = 2008
    86 = 2017
    87 = 2011
    88 = 2018
    89 = 2016
    90 = 2019
    91 = 2017
    92 = 2015
    93 = 2016
    94 = 2012
    95 = 2015
    96 = 2021
    97 = 2015
    98 = 2020
    99 = 2012
    100 = 2015
    101 = 2018
    102 = 2021
    103 = 2020
    104 = 2018
    105 = 2019
    106 = 2014
    107 = 2019
    108 = 2018
    109 = 2017
    110 = 2020
}

foreach ($row in $yearData.Keys) {
    $ws.Cells.Item([int]$row, 5).Value = $yearData[$row]
}

# Restore the scroll position / active selection recorded in the sheet view.
$ws.Range("F110").Select()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
